$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

$ws1.Range("F5").Value = 362
$ws1.Range("F6").Value = 164
$ws1.Range("F7").Value = 157
$ws1.Range("F8").Value = 765
$ws1.Range("F9").Value = 4151
$ws1.Range("F12").Value = 169
$ws1.Range("F14").Value = 5948
$ws1.Range("F15").Value = 461
$ws1.Range("F16").Value = 2294
$ws1.Range("F18").Value = 162
$ws1.Range("F19").Value = 453
$ws1.Range("F20").Value = 8952
$ws1.Range("C22").Value = "北京·喘气动漫嘉年华·暑期狂欢"
$ws1.Range("F22").Value = 2123
$ws1.Range("F24").Value = 2283
$ws1.Range("F25").Value = 2391
$ws1.Range("F26").Value = 1375
$ws1.Range("F27").Value = 224
$ws1.Range("F28").Value = 1934
$ws1.Range("F29").Value = 40
$ws1.Range("F30").Value = 54
$ws1.Range("F35").Value = 32
$ws1.Range("F37").Value = 1220
$ws1.Range("F41").Value = 226
$ws1.Range("F42").Value = 1501
$ws1.Range("F43").Value = 2431
$ws1.Range("F45").Value = 905
$ws1.Range("F46").Value = 284
$ws1.Range("F48").Value = 3
$ws2.Range("C11").Value = "北京·“记忆重启”—— 《鬼怪》《请回答1988》《来自星星的你》一起追过的影视剧名曲音乐会（取消）"
$ws2.Range("G11").Value = "不可售"
$ws2.Range("F22").Value = 26
$ws2.Range("F23").Value = 26
$ws3.Range("F2").Value = 683
$ws3.Range("F3").Value = 880
$ws4.Range("F3").Value = 683
$ws4.Range("F4").Value = 880
$ws4.Range("F6").Value = 362
$ws4.Range("F8").Value = 164
$ws4.Range("F11").Value = 157
$ws4.Range("F12").Value = 765
$ws4.Range("F13").Value = 4151
$ws4.Range("F14").Value = 4151
$ws4.Range("F16").Value = 169
$ws4.Range("C18").Value = "北京·广播剧《蝉女》专场活动"
$ws4.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=86947"
$ws4.Range("G18").Value = 288
$ws4.Range("E18").Value = "2024.08.10 11:50-08.10 15:10"
$ws4.Range("I18").Value = "//i2.hdslb.com/bfs/openplatform/202406/ycrRjEPg1718176423186.jpeg"
$ws4.Range("F18").Value = 95
$ws4.Range("C19").Value = "北京·梦次元动漫展M30"
$ws4.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=83828"
$ws4.Range("G19").Value = 80
$ws4.Range("E19").Value = "2024.08.10 10:00-08.11 17:00"
$ws4.Range("I19").Value = "//i1.hdslb.com/bfs/openplatform/202405/Qr2Bd5W41715931423636.jpeg"
$ws4.Range("F19").Value = 5948
$ws4.Range("G20").Value = 688
$ws4.Range("D20").Value = "北京展览馆 北京展览馆"
$ws4.Range("F20").Value = 461
$ws4.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=87199"
$ws4.Range("C20").Value = "北京·知名演员 捣宫姬奈&蕨野友也 专场活动"
$ws4.Range("I20").Value = "//i0.hdslb.com/bfs/openplatform/202406/6MjCfuLT1718175794197.png"
$ws4.Range("E20").Value = "2024.08.10 10:30-08.10 15:50"
$ws4.Range("I21").Value = "//i1.hdslb.com/bfs/openplatform/202405/4jQBoo241716968548735.jpeg"
$ws4.Range("B21").Value = "2024-08-10"
$ws4.Range("E21").Value = "2024.08.10 10:00-08.10 17:00"
$ws4.Range("D21").Value = "永外高庄138号 北京大红门国际会展中心"
$ws4.Range("F21").Value = 2294
$ws4.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=86590"
$ws4.Range("G21").Value = 70
$ws4.Range("C21").Value = "北京·第五人格ONLY2.0"
$ws4.Range("F22").Value = 162
$ws4.Range("F23").Value = 453
$ws4.Range("F24").Value = 8952
$ws4.Range("C27").Value = "北京·喘气动漫嘉年华·暑期狂欢"
$ws4.Range("F27").Value = 2124
$ws4.Range("F28").Value = 2283
$ws4.Range("F29").Value = 2391
$ws4.Range("F30").Value = 1375
$ws4.Range("F31").Value = 224
$ws4.Range("F32").Value = 1934
$ws4.Range("F33").Value = 40
$ws4.Range("F34").Value = 54
$ws4.Range("F37").Value = 32
$ws4.Range("F39").Value = 1220
$ws4.Range("F42").Value = 226
$ws4.Range("F43").Value = 1501
$ws4.Range("F44").Value = 2431
$ws4.Range("F45").Value = 905
$ws4.Range("F46").Value = 284
$ws4.Range("F51").Value = 26
